# Corrected Q2 Done Q3 Done DF resultat filled Done Choice 2
#
# - Remove rows 12-18 (extra DF records that shouldn't be in the result).
# - Remove column F (birthPlace) entirely.
# - Fill the previously-empty birthDate cells (E3, E5, E6, E8, E9) with "nan"
#   to reflect the filled DF result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 12 through 18.
$ws.Range("A12:A18").EntireRow.Delete() | Out-Null

# Remove column F (birthPlace).
$ws.Columns.Item(6).Delete() | Out-Null

# Fill the empty birthDate values with the string "nan".
$ws.Range("E3").Value = "nan"
$ws.Range("E5").Value = "nan"
$ws.Range("E6").Value = "nan"
$ws.Range("E8").Value = "nan"
$ws.Range("E9").Value = "nan"
